$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("login")
$ws.Range("D9").Value = "NA"
